{"js": "// Table S5 (DDM a) \u2014 recalculate and change optics:\n// update the p-values, a beta estimate, and two 95% CI cells that\n// changed after the model was recalculated.\nconst replacements = [\n  [\".367\", \".362\"],               // (Intercept) row, p column\n  [\".266\", \".239\"],               // CLhigh row, p column\n  [\"[-0.29, 0.08]\", \"[-0.30, 0.08]\"], // CLhigh row, 95% CI column\n  [\"-0.17\", \"-0.16\"],             // Eval1 row, b column\n  [\".082\", \".091\"],               // Eval1 row, p column\n  [\"[-0.36, 0.02]\", \"[-0.35, 0.03]\"], // Eval1 row, 95% CI column\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + findText + \"', found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Table S5 (DDM a) \u2014 recalculate and change optics:\n# update the p-values, a beta estimate, and two 95% CI cells that\n# changed after the model was recalculated.\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute([ref]$findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Replace-Once: could not find text '$findText'\"\n    }\n}\n\nReplace-Once \".367\" \".362\"                         # (Intercept) row, p column\nReplace-Once \".266\" \".239\"                         # CLhigh row, p column\nReplace-Once \"[-0.29, 0.08]\" \"[-0.30, 0.08]\"       # CLhigh row, 95% CI column\nReplace-Once \"-0.17\" \"-0.16\"                       # Eval1 row, b column\nReplace-Once \".082\" \".091\"                         # Eval1 row, p column\nReplace-Once \"[-0.36, 0.02]\" \"[-0.35, 0.03]\"       # Eval1 row, 95% CI column\n"}
